# Applies the "updates and additions user story" edit:
#   1. Splits the "In this project..." paragraph's single run into three
#      runs, inserting " and Contact Form" after "To Do List App".
#   2. Adds a brand-new paragraph right after the
#      "How to Use the App/What Functionality it Offers:" paragraph that
#      describes the app's/contact form's functionality, including the
#      two <w:proofErr.../> grammar-check markers from the original edit.
#
# Both edits are done via Range.InsertXML() with a full WordprocessingML
# "pkg:package" fragment so the exact run/markup boundaries from the
# target OOXML are reproduced verbatim (no incidental rPr/formatting is
# introduced, and standalone markup like w:proofErr survives untouched).

$d = $word.ActiveDocument

function Set-RangeOOXML($range, [string]$bodyXml) {
    $pkg = '<?xml version="1.0" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body>' + $bodyXml + '</w:body>' +
        '</w:document>' +
        '</pkg:xmlData>' +
        '</pkg:part>' +
        '</pkg:package>'
    $range.InsertXML($pkg)
}

# --- 1. "In this project..." -> split into 3 runs ------------------------

$target = "In this project we were asked to design a To Do List App. In the app, the user should be able to add list items, remove/make the item disappear from the list, and mark list items as complete. In addition, we needed to make a page that allows a user to complete a contact form. In the form they should be able to submit their information and a message/comment. "

$introPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13) -eq $target) {
        $introPara = $p
        break
    }
}

if ($introPara -eq $null) {
    throw "Could not locate the intro paragraph to split."
}

$introRange = $d.Range($introPara.Range.Start, $introPara.Range.End)
$introBody = '<w:p>' +
    '<w:r><w:t>In this project we were asked to design a To Do List App</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> and Contact Form</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">. In the app, the user should be able to add list items, remove/make the item disappear from the list, and mark list items as complete. In addition, we needed to make a page that allows a user to complete a contact form. In the form they should be able to submit their information and a message/comment. </w:t></w:r>' +
    '</w:p>'
Set-RangeOOXML $introRange $introBody

# --- 2. New "How to use" paragraph with proofErr markers ------------------

$howToPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13) -eq "How to Use the App/What Functionality it Offers:") {
        $howToPara = $p
        break
    }
}

if ($howToPara -eq $null) {
    throw "Could not locate the 'How to Use the App' paragraph."
}

# Insert a brand-new empty paragraph right after it, then fill it with the
# user-story text (including the grammar-checker proofErr markers) via the
# same InsertXML technique.
$howToPara.Range.InsertParagraphAfter()
$newPara = $howToPara.Next()

$newParaRange = $d.Range($newPara.Range.Start, $newPara.Range.End)
$newBody = '<w:p>' +
    '<w:r><w:t xml:space="preserve">The app is simple. The user can type in their list items. If they want to change </w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:t>them</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:t xml:space="preserve"> they can click the icons &#x201c;X&#x201d; for delete and &#x201c;pencil and paper&#x201d; to edit the item. When they wish to update the </w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:t>item</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:t xml:space="preserve"> they will be able to type in the item they wish to add. The user can also click the item and it will get darker to show that this item has been completed. The contact form is simple, they just click on the Contact form from the navigation bar and fill in the information below. Once they have filled in the information, the user can submit.</w:t></w:r>' +
    '</w:p>'
Set-RangeOOXML $newParaRange $newBody

Write-Output "Done."
